$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "UNIPROT_ID"
$ws.Range("D1").Value = "GENBANK_ID"
$ws.Range("F7").Select() | Out-Null
